$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-05-08 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-09 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("82×34=2788", $true, $false, $false, $false, $false, $true, 1, $false, "15×96=1440", 2) | Out-Null
$d.Content.Find.Execute("29×70=2030", $true, $false, $false, $false, $false, $true, 1, $false, "82×45=3690", 2) | Out-Null
$d.Content.Find.Execute("78×53=4134", $true, $false, $false, $false, $false, $true, 1, $false, "53×24=1272", 2) | Out-Null
$d.Content.Find.Execute("46×32=1472", $true, $false, $false, $false, $false, $true, 1, $false, "63×73=4599", 2) | Out-Null
$d.Content.Find.Execute("41×100=4100", $true, $false, $false, $false, $false, $true, 1, $false, "63×81=5103", 2) | Out-Null
$d.Content.Find.Execute("39×73=2847", $true, $false, $false, $false, $false, $true, 1, $false, "76×27=2052", 2) | Out-Null
$d.Content.Find.Execute("26×72=1872", $true, $false, $false, $false, $false, $true, 1, $false, "30×97=2910", 2) | Out-Null
$d.Content.Find.Execute("86×89=7654", $true, $false, $false, $false, $false, $true, 1, $false, "92×26=2392", 2) | Out-Null
$d.Content.Find.Execute("70×71=4970", $true, $false, $false, $false, $false, $true, 1, $false, "16×100=1600", 2) | Out-Null
$d.Content.Find.Execute("62×58=3596", $true, $false, $false, $false, $false, $true, 1, $false, "87×76=6612", 2) | Out-Null
$d.Content.Find.Execute("70×36=2520", $true, $false, $false, $false, $false, $true, 1, $false, "32×46=1472", 2) | Out-Null
$d.Content.Find.Execute("90×46=4140", $true, $false, $false, $false, $false, $true, 1, $false, "89×24=2136", 2) | Out-Null
$d.Content.Find.Execute("25×93=2325", $true, $false, $false, $false, $false, $true, 1, $false, "99×74=7326", 2) | Out-Null
$d.Content.Find.Execute("19×58=1102", $true, $false, $false, $false, $false, $true, 1, $false, "47×15=705", 2) | Out-Null
$d.Content.Find.Execute("40×95=3800", $true, $false, $false, $false, $false, $true, 1, $false, "78×76=5928", 2) | Out-Null
$d.Content.Find.Execute("39×36=1404", $true, $false, $false, $false, $false, $true, 1, $false, "35×81=2835", 2) | Out-Null
$d.Content.Find.Execute("10×80=800", $true, $false, $false, $false, $false, $true, 1, $false, "28×95=2660", 2) | Out-Null
$d.Content.Find.Execute("98×81=7938", $true, $false, $false, $false, $false, $true, 1, $false, "44×57=2508", 2) | Out-Null
$d.Content.Find.Execute("18×79=1422", $true, $false, $false, $false, $false, $true, 1, $false, "41×83=3403", 2) | Out-Null
$d.Content.Find.Execute("65×50=3250", $true, $false, $false, $false, $false, $true, 1, $false, "81×41=3321", 2) | Out-Null
$d.Content.Find.Execute("26×51=1326", $true, $false, $false, $false, $false, $true, 1, $false, "34×27=918", 2) | Out-Null
$d.Content.Find.Execute("11×72=792", $true, $false, $false, $false, $false, $true, 1, $false, "50×36=1800", 2) | Out-Null
$d.Content.Find.Execute("29×68=1972", $true, $false, $false, $false, $false, $true, 1, $false, "41×61=2501", 2) | Out-Null
$d.Content.Find.Execute("92×77=7084", $true, $false, $false, $false, $false, $true, 1, $false, "54×28=1512", 2) | Out-Null
$d.Content.Find.Execute("60×77=4620", $true, $false, $false, $false, $false, $true, 1, $false, "65×85=5525", 2) | Out-Null
$d.Content.Find.Execute("12×83=996", $true, $false, $false, $false, $false, $true, 1, $false, "49×33=1617", 2) | Out-Null
$d.Content.Find.Execute("29×33=957", $true, $false, $false, $false, $false, $true, 1, $false, "24×35=840", 2) | Out-Null
$d.Content.Find.Execute("51×76=3876", $true, $false, $false, $false, $false, $true, 1, $false, "38×59=2242", 2) | Out-Null
$d.Content.Find.Execute("60×12=720", $true, $false, $false, $false, $false, $true, 1, $false, "33×57=1881", 2) | Out-Null
$d.Content.Find.Execute("38×88=3344", $true, $false, $false, $false, $false, $true, 1, $false, "16×57=912", 2) | Out-Null
$d.Content.Find.Execute("88×59=5192", $true, $false, $false, $false, $false, $true, 1, $false, "24×25=600", 2) | Out-Null
$d.Content.Find.Execute("40×99=3960", $true, $false, $false, $false, $false, $true, 1, $false, "23×81=1863", 2) | Out-Null
$d.Content.Find.Execute("72×77=5544", $true, $false, $false, $false, $false, $true, 1, $false, "12×43=516", 2) | Out-Null
$d.Content.Find.Execute("32×23=736", $true, $false, $false, $false, $false, $true, 1, $false, "57×51=2907", 2) | Out-Null
$d.Content.Find.Execute("17×70=1190", $true, $false, $false, $false, $false, $true, 1, $false, "88×93=8184", 2) | Out-Null
$d.Content.Find.Execute("78×54=4212", $true, $false, $false, $false, $false, $true, 1, $false, "50×89=4450", 2) | Out-Null
$d.Content.Find.Execute("67×71=4757", $true, $false, $false, $false, $false, $true, 1, $false, "17×82=1394", 2) | Out-Null
$d.Content.Find.Execute("55×27=1485", $true, $false, $false, $false, $false, $true, 1, $false, "17×17=289", 2) | Out-Null
$d.Content.Find.Execute("63×35=2205", $true, $false, $false, $false, $false, $true, 1, $false, "61×73=4453", 2) | Out-Null
$d.Content.Find.Execute("62×100=6200", $true, $false, $false, $false, $false, $true, 1, $false, "60×83=4980", 2) | Out-Null
$d.Content.Find.Execute("97×42=4074", $true, $false, $false, $false, $false, $true, 1, $false, "73×11=803", 2) | Out-Null
$d.Content.Find.Execute("50×30=1500", $true, $false, $false, $false, $false, $true, 1, $false, "88×70=6160", 2) | Out-Null
$d.Content.Find.Execute("83×34=2822", $true, $false, $false, $false, $false, $true, 1, $false, "89×26=2314", 2) | Out-Null
$d.Content.Find.Execute("61×74=4514", $true, $false, $false, $false, $false, $true, 1, $false, "48×17=816", 2) | Out-Null
$d.Content.Find.Execute("100×15=1500", $true, $false, $false, $false, $false, $true, 1, $false, "93×73=6789", 2) | Out-Null
$d.Content.Find.Execute("83×80=6640", $true, $false, $false, $false, $false, $true, 1, $false, "11×97=1067", 2) | Out-Null
$d.Content.Find.Execute("23×27=621", $true, $false, $false, $false, $false, $true, 1, $false, "10×11=110", 2) | Out-Null
$d.Content.Find.Execute("42×45=1890", $true, $false, $false, $false, $false, $true, 1, $false, "79×46=3634", 2) | Out-Null
$d.Content.Find.Execute("63×64=4032", $true, $false, $false, $false, $false, $true, 1, $false, "24×86=2064", 2) | Out-Null
$d.Content.Find.Execute("37×74=2738", $true, $false, $false, $false, $false, $true, 1, $false, "56×49=2744", 2) | Out-Null
$d.Content.Find.Execute("49×32=1568", $true, $false, $false, $false, $false, $true, 1, $false, "58×88=5104", 2) | Out-Null
$d.Content.Find.Execute("25×80=2000", $true, $false, $false, $false, $false, $true, 1, $false, "64×85=5440", 2) | Out-Null
$d.Content.Find.Execute("88×95=8360", $true, $false, $false, $false, $false, $true, 1, $false, "88×32=2816", 2) | Out-Null
$d.Content.Find.Execute("53×49=2597", $true, $false, $false, $false, $false, $true, 1, $false, "44×59=2596", 2) | Out-Null
$d.Content.Find.Execute("57×85=4845", $true, $false, $false, $false, $false, $true, 1, $false, "86×22=1892", 2) | Out-Null
$d.Content.Find.Execute("74×60=4440", $true, $false, $false, $false, $false, $true, 1, $false, "29×39=1131", 2) | Out-Null
$d.Content.Find.Execute("94×21=1974", $true, $false, $false, $false, $false, $true, 1, $false, "94×43=4042", 2) | Out-Null
$d.Content.Find.Execute("10×95=950", $true, $false, $false, $false, $false, $true, 1, $false, "70×94=6580", 2) | Out-Null
$d.Content.Find.Execute("52×87=4524", $true, $false, $false, $false, $false, $true, 1, $false, "28×80=2240", 2) | Out-Null
$d.Content.Find.Execute("76×96=7296", $true, $false, $false, $false, $false, $true, 1, $false, "85×28=2380", 2) | Out-Null
$d.Content.Find.Execute("32×58=1856", $true, $false, $false, $false, $false, $true, 1, $false, "56×70=3920", 2) | Out-Null
$d.Content.Find.Execute("66×59=3894", $true, $false, $false, $false, $false, $true, 1, $false, "70×84=5880", 2) | Out-Null
$d.Content.Find.Execute("84×83=6972", $true, $false, $false, $false, $false, $true, 1, $false, "27×82=2214", 2) | Out-Null
$d.Content.Find.Execute("33×87=2871", $true, $false, $false, $false, $false, $true, 1, $false, "66×93=6138", 2) | Out-Null
$d.Content.Find.Execute("43×53=2279", $true, $false, $false, $false, $false, $true, 1, $false, "24×26=624", 2) | Out-Null
$d.Content.Find.Execute("88×91=8008", $true, $false, $false, $false, $false, $true, 1, $false, "53×69=3657", 2) | Out-Null
$d.Content.Find.Execute("83×32=2656", $true, $false, $false, $false, $false, $true, 1, $false, "59×10=590", 2) | Out-Null
$d.Content.Find.Execute("80×87=6960", $true, $false, $false, $false, $false, $true, 1, $false, "70×45=3150", 2) | Out-Null
$d.Content.Find.Execute("52×31=1612", $true, $false, $false, $false, $false, $true, 1, $false, "83×98=8134", 2) | Out-Null
$d.Content.Find.Execute("17×64=1088", $true, $false, $false, $false, $false, $true, 1, $false, "92×34=3128", 2) | Out-Null
$d.Content.Find.Execute("76×38=2888", $true, $false, $false, $false, $false, $true, 1, $false, "40×49=1960", 2) | Out-Null
$d.Content.Find.Execute("53×14=742", $true, $false, $false, $false, $false, $true, 1, $false, "22×23=506", 2) | Out-Null
$d.Content.Find.Execute("47×56=2632", $true, $false, $false, $false, $false, $true, 1, $false, "56×25=1400", 2) | Out-Null
$d.Content.Find.Execute("31×74=2294", $true, $false, $false, $false, $false, $true, 1, $false, "87×12=1044", 2) | Out-Null
$d.Content.Find.Execute("15×77=1155", $true, $false, $false, $false, $false, $true, 1, $false, "89×44=3916", 2) | Out-Null
$d.Content.Find.Execute("96×55=5280", $true, $false, $false, $false, $false, $true, 1, $false, "14×45=630", 2) | Out-Null
$d.Content.Find.Execute("49×37=1813", $true, $false, $false, $false, $false, $true, 1, $false, "14×10=140", 2) | Out-Null
$d.Content.Find.Execute("42×94=3948", $true, $false, $false, $false, $false, $true, 1, $false, "14×74=1036", 2) | Out-Null
$d.Content.Find.Execute("88×75=6600", $true, $false, $false, $false, $false, $true, 1, $false, "88×62=5456", 2) | Out-Null
$d.Content.Find.Execute("67×41=2747", $true, $false, $false, $false, $false, $true, 1, $false, "48×14=672", 2) | Out-Null
$d.Content.Find.Execute("35×39=1365", $true, $false, $false, $false, $false, $true, 1, $false, "92×18=1656", 2) | Out-Null
$d.Content.Find.Execute("89×64=5696", $true, $false, $false, $false, $false, $true, 1, $false, "41×70=2870", 2) | Out-Null
$d.Content.Find.Execute("43×55=2365", $true, $false, $false, $false, $false, $true, 1, $false, "84×87=7308", 2) | Out-Null
$d.Content.Find.Execute("14×48=672", $true, $false, $false, $false, $false, $true, 1, $false, "73×52=3796", 2) | Out-Null
$d.Content.Find.Execute("13×43=559", $true, $false, $false, $false, $false, $true, 1, $false, "33×88=2904", 2) | Out-Null
$d.Content.Find.Execute("35×33=1155", $true, $false, $false, $false, $false, $true, 1, $false, "97×29=2813", 2) | Out-Null
$d.Content.Find.Execute("82×84=6888", $true, $false, $false, $false, $false, $true, 1, $false, "10×76=760", 2) | Out-Null
$d.Content.Find.Execute("30×23=690", $true, $false, $false, $false, $false, $true, 1, $false, "30×66=1980", 2) | Out-Null
$d.Content.Find.Execute("85×81=6885", $true, $false, $false, $false, $false, $true, 1, $false, "47×97=4559", 2) | Out-Null
$d.Content.Find.Execute("22×72=1584", $true, $false, $false, $false, $false, $true, 1, $false, "49×74=3626", 2) | Out-Null
$d.Content.Find.Execute("46×19=874", $true, $false, $false, $false, $false, $true, 1, $false, "36×81=2916", 2) | Out-Null
$d.Content.Find.Execute("81×49=3969", $true, $false, $false, $false, $false, $true, 1, $false, "10×63=630", 2) | Out-Null
$d.Content.Find.Execute("83×65=5395", $true, $false, $false, $false, $false, $true, 1, $false, "95×57=5415", 2) | Out-Null
$d.Content.Find.Execute("40×50=2000", $true, $false, $false, $false, $false, $true, 1, $false, "28×32=896", 2) | Out-Null
$d.Content.Find.Execute("40×33=1320", $true, $false, $false, $false, $false, $true, 1, $false, "98×50=4900", 2) | Out-Null
$d.Content.Find.Execute("68×21=1428", $true, $false, $false, $false, $false, $true, 1, $false, "47×34=1598", 2) | Out-Null
$d.Content.Find.Execute("72×99=7128", $true, $false, $false, $false, $false, $true, 1, $false, "38×47=1786", 2) | Out-Null
$d.Content.Find.Execute("47×39=1833", $true, $false, $false, $false, $false, $true, 1, $false, "93×39=3627", 2) | Out-Null
$d.Content.Find.Execute("97×44=4268", $true, $false, $false, $false, $false, $true, 1, $false, "35×74=2590", 2) | Out-Null
$d.Content.Find.Execute("74×83=6142", $true, $false, $false, $false, $false, $true, 1, $false, "56×21=1176", 2) | Out-Null
